$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $text) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.ClearFormats()
}

# Row 2: Flamengo unchanged, possession % updated
Set-TextValue $ws.Range("B2") "59.2%"

# Row 3: Corinthians unchanged, possession % updated
Set-TextValue $ws.Range("B3") "55.8%"

# Row 4: was Fluminense/53.6% -> now Bahia/53.9%
$ws.Range("A4").Value = "Bahia"
Set-TextValue $ws.Range("B4") "53.9%"

# Row 5: was Bahia/53.5% -> now Fluminense/53.6%
$ws.Range("A5").Value = "Fluminense"
Set-TextValue $ws.Range("B5") "53.6%"

# Row 6: was Atletico MG/53.3% -> now Vasco da Gama/53.5%
$ws.Range("A6").Value = "Vasco da Gama"
Set-TextValue $ws.Range("B6") "53.5%"
